$wb = $excel.ActiveWorkbook

function Get-HyperlinkByCell($worksheet, $cellAddr) {
    foreach ($hl in $worksheet.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellAddr) {
            return $hl
        }
    }
    return $null
}

function Add-MatchingHyperlink($worksheet, $targetCell, $sourceCellAddr, $displayText) {
    $srcHl = Get-HyperlinkByCell $worksheet $sourceCellAddr
    $rng = $worksheet.Range($targetCell)
    $rng.Value = $displayText
    $worksheet.Hyperlinks.Add($rng, $srcHl.Address, "", "", $displayText) | Out-Null
    $rng.Font.Underline = 2
    $rng.Font.Color = 15570276
}

# ---------------------------------------------------------------------------
# 1. Update the "Status" column everywhere it says "Ready for handoff" so it
#    now reads "Handed back: in sync with en-US" (Overview + zh-cn + de-de).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: fill in "Latest Target File" (F) and "Latest Handback File"
#    (G) columns, pointing at the same targets as the Source File (A) and
#    Latest Handoff File (D) columns, and stamp the real handback datetime
#    into column H (was the zero-date placeholder).
# ---------------------------------------------------------------------------
Add-MatchingHyperlink $wsZhCn "F2" "`$A`$2" "7b09e590-e9f9-4ec8-82b0-57ecea581437.md"
Add-MatchingHyperlink $wsZhCn "G2" "`$D`$2" "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.zh-cn.xlf"
Add-MatchingHyperlink $wsZhCn "F3" "`$A`$3" "e6081861-f56e-4a77-9236-9901c9cdc7e9.md"
Add-MatchingHyperlink $wsZhCn "G3" "`$D`$3" "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.zh-cn.xlf"

$wsZhCn.Range("H2").Value = "2016-03-25 09:01:39"
$wsZhCn.Range("H3").Value = "2016-03-25 09:01:39"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same treatment, using the de-de targets / timestamp.
# ---------------------------------------------------------------------------
Add-MatchingHyperlink $wsDeDe "F2" "`$A`$2" "7b09e590-e9f9-4ec8-82b0-57ecea581437.md"
Add-MatchingHyperlink $wsDeDe "G2" "`$D`$2" "7b09e590-e9f9-4ec8-82b0-57ecea581437.4fdc4ad2453416e45d2658ced24248438eeb5397.de-de.xlf"
Add-MatchingHyperlink $wsDeDe "F3" "`$A`$3" "e6081861-f56e-4a77-9236-9901c9cdc7e9.md"
Add-MatchingHyperlink $wsDeDe "G3" "`$D`$3" "e6081861-f56e-4a77-9236-9901c9cdc7e9.65b374d3eee2adffd0001118f274db0000a84dd7.de-de.xlf"

$wsDeDe.Range("H2").Value = "2016-03-25 09:01:59"
$wsDeDe.Range("H3").Value = "2016-03-25 09:01:59"
